# Changes of 11th July 2022
# Updates a batch of tracking-number cells (ShipmentTrackNum / PackageTrackNum)
# in Sheet1 with newer tracking numbers, and flips the Q3 pass/fail flag.
#
# We write each value through a scratch cell as a `="text"` formula and then
# PasteSpecial(xlPasteValues) it into the target cell. Doing this (instead of
# a plain `Range.Value = "<digits>"` assignment) is required because these
# values are long purely-numeric strings: a direct `.Value =` assignment
# auto-coerces them to a numeric cell (losing the shared-string/text type and
# silently attaching a new cell style), whereas copy/pasting the *result* of
# a text-producing formula keeps the cell's original style untouched and
# stores the value as shared-string text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $text)

    $scratch = $ws.Range("ZZ9000")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()

    $cell = $ws.Range($addr)
    $cell.PasteSpecial(-4163)  # xlPasteValues

    $scratch.Clear()
}

Set-TextValue $ws "C2"  "320018675783"
Set-TextValue $ws "C3"  "320018675794"
Set-TextValue $ws "Q3"  "FAIL"
Set-TextValue $ws "C4"  "320018675820"
Set-TextValue $ws "C5"  "320018675842"
Set-TextValue $ws "D5"  "320018675842"
Set-TextValue $ws "C6"  "320018675886"
Set-TextValue $ws "D6"  "320018675886"
Set-TextValue $ws "C7"  "320018675901"
Set-TextValue $ws "D7"  "320018675901"
Set-TextValue $ws "C8"  "320018675934"
Set-TextValue $ws "C9"  "320018675956"
Set-TextValue $ws "C10" "320018675989"
Set-TextValue $ws "C11" "320018676025"
Set-TextValue $ws "C12" "320018676069"
Set-TextValue $ws "C13" "320018676080"
Set-TextValue $ws "D13" "320018676080"
Set-TextValue $ws "C14" "320018676117"
Set-TextValue $ws "D14" "320018676117"
Set-TextValue $ws "C15" "320018676139"
Set-TextValue $ws "D15" "320018676139"
Set-TextValue $ws "C16" "320018676161"
Set-TextValue $ws "D16" "320018676161"
Set-TextValue $ws "C17" "320018676183"
Set-TextValue $ws "D17" "320018676183"
Set-TextValue $ws "C18" "320018677775"
Set-TextValue $ws "C19" "320018677797"
Set-TextValue $ws "C20" "320018677948"
Set-TextValue $ws "C21" "320018677992"
Set-TextValue $ws "C22" "320018678028"

Write-Host "Applied tracking-number + Q3 status updates."
